$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.893.00'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.551.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('E8').Value = '  +0.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.50'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0856'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.772.77'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.549.76'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.900.81'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.61'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '213.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.15'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.78'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('E26').Value = '  +2.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.85'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0458'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('E32').Value = '  +1.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.370.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.98%  '
$ws.Range('E34').Value = '  +1.89%  '
$ws.Range('E35').Value = '  +3.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.971'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.42%  '
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('E38').Value = '  +1.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.522'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.983'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +3.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.61'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.18%  '
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.686.21'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0506'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('E50').Value = '  +1.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.30%  '
